# Apply the GHG GCHeadspace template fixes described in the commit:
# "New fixes to GHG processing templates"
#
# Net semantic change on Sheet1:
#   - Header O1: "Water_Temp_DegC" -> "WaterT_C"
#   - Header P1: "Ambient_air"      -> "Air_Location"
#   - New row 2: O2 and P2 both get the note "Has been changed thru samplings"
#   - Column O widened to fit the new (longer) note text
#   - Active selection moved to O9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two headers (values only; existing cell styles/formats stay) ---
$ws.Cells.Item(1, 15).Value = "WaterT_C"       # O1
$ws.Cells.Item(1, 16).Value = "Air_Location"   # P1

# --- Add the new second row of notes under the renamed headers ---
$ws.Cells.Item(2, 15).Value = "Has been changed thru samplings"  # O2
$ws.Cells.Item(2, 16).Value = "Has been changed thru samplings"  # P2

# --- Widen column O so the longer text fits (matches the new bestFit width) ---
$ws.Columns.Item(15).AutoFit()
$ws.Columns.Item(15).ColumnWidth = 27.33

# --- Update the saved selection/active cell ---
$null = $ws.Range("O9").Select()
